# New weekly price observation for "Haba" arrives and is inserted as a new
# data row right after the header/second rows, pushing the existing rows
# (previously 30..138) down by one (now 31..139).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row before the current row 30 - this shifts rows
# 30..138 down to 31..139, preserving all of their existing data/formatting.
$ws.Rows(30).Insert()

# Populate the newly inserted row 30 with the new observation.
$ws.Cells.Item(30, 1).Value  = 3
$ws.Cells.Item(30, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(30, 3).Value  = "Coquimbo"
$ws.Cells.Item(30, 4).Value  = 44707
$ws.Cells.Item(30, 5).Value  = 5
$ws.Cells.Item(30, 6).Value  = 100112026
$ws.Cells.Item(30, 7).Value  = "Haba"
$ws.Cells.Item(30, 8).Value  = "Sin especificar"
$ws.Cells.Item(30, 9).Value  = "Primera"
$ws.Cells.Item(30, 10).Value = 85
$ws.Cells.Item(30, 11).Value = 20000
$ws.Cells.Item(30, 12).Value = 21000
$ws.Cells.Item(30, 13).Value = 20529
$ws.Cells.Item(30, 14).Value = "`$/malla 25 kilos"
$ws.Cells.Item(30, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(30, 16).Value = 821
$ws.Cells.Item(30, 17).Value = 25
$ws.Cells.Item(30, 18).Value = "Hortaliza"
